# Add explanation to Leetcode - 131. Palindrome Partitioning - Backtracking
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

$category = "Backtracking"
$difficulty = "Medium"
$name = "131. Palindrome Partitioning"
$notes = "The basic idea is - Whats the biggest palindrome partition we can get from i to j`r`nWe run a for loop j in range(i, len(s)), in that we check for partition slice from i to j, if its valid append that slice to curr and run dfs(j + 1) to solve for the rest of the array, after that pop the slice we appended, so that we can look at a bigger slice in the next iteration of the loop.`r`nIn the next iteration j will move and our window will get bigger, if the slice is a partiton append it to curr again and look at the  rest of the array, otherwise we dont need to look at this branch any longer.`r`nBASECASE: We need to find the final array after finding out all combos of the partitions, so add to res only when i >= n"

$row = 24

$ws.Cells.Item($row, 1).Value = $category
$ws.Cells.Item($row, 2).Value = $difficulty
$ws.Cells.Item($row, 3).Value = $name
$ws.Cells.Item($row, 4).Value = $notes

$ws.Cells.Item($row, 2).Style = "Neutral"
$ws.Cells.Item($row, 3).Style = "Neutral"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 3), "https://leetcode.com/problems/palindrome-partitioning/") | Out-Null

# Re-assert the Neutral style since adding the hyperlink re-applies the default Hyperlink style
$ws.Cells.Item($row, 3).Style = "Neutral"

$ws.Rows.Item($row).RowHeight = 87

$ws.Range("A25").Select() | Out-Null
